# Add a new minutes sheet "1603" (16 March) positioned after "1103",
# based on a copy of "1103", with the Action Review / New Actions /
# Decisions tables updated for the new meeting.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the "1103" sheet, place the copy right after it, rename.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("1103")
$template.Copy($null, $template)
$ws = $wb.Worksheets.Item($template.Index + 1)
$ws.Name = "1603"

# Donor sheet used to borrow cell formats (styles) that already exist
# elsewhere in the workbook but are not present (in the needed spot) on
# the freshly duplicated sheet.
$donor = $wb.Worksheets.Item("0403")

# ---------------------------------------------------------------------
# 2. Restructure rows so the table sizes match the new minutes:
#      Action Review data:      2 rows (6:7)   -> 5 rows (6:10)
#      New Actions/Owners data: 5 rows (11:15) -> 3 rows (14:16, after shift)
#      Decisions data:          2 rows (18:19) -> 3 rows (19:21, after shift)
# ---------------------------------------------------------------------
$ws.Rows("8:10").Insert()     # grow Action Review data block by 3 rows
$ws.Rows("17:18").Delete()    # shrink New Actions/Owners data block by 2 rows
$ws.Rows("21:21").Insert()    # grow Decisions data block by 1 row

# ---------------------------------------------------------------------
# 3. Meeting date (row 1)
# ---------------------------------------------------------------------
$ws.Range("B1").Value = 44271

# ---------------------------------------------------------------------
# 4. Action Review table (rows 6:10) - fix styles, then values.
# ---------------------------------------------------------------------
$donor.Range("A11").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$donor.Range("B11").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$donor.Range("C11").Copy()
$ws.Range("C6").PasteSpecial(-4122)

$donor.Range("A11:C15").Copy()
$ws.Range("A7:C10").PasteSpecial(-4122)

$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30

$ws.Range("A6").Value = "Finalise prototype"
$ws.Range("B6").Value = "Hamza"
$ws.Range("C6").Value = "Finalising"

$ws.Range("A7").Value = "Finalise survey monkey"
$ws.Range("B7").Value = "Emily"
$ws.Range("C7").Value = "Completed"

$ws.Range("A8").Value = "Develop data model"
$ws.Range("B8").Value = "Jess"
$ws.Range("C8").Value = "Basic version completed"

$ws.Range("A9").Value = "Write project outline for github"
$ws.Range("B9").Value = "Jess"
$ws.Range("C9").Value = "Completed"

$ws.Range("A10").Value = "Create header page for website"
$ws.Range("B10").Value = "Stanni"
$ws.Range("C10").Value = "Completed"

# ---------------------------------------------------------------------
# 5. New Actions and Owners table (rows 14:16) - fix styles, then values.
# ---------------------------------------------------------------------
$donor.Range("A11").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$donor.Range("B11").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$donor.Range("C11").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 30

$ws.Range("A14").Value = "Share survey with friends and family"
$ws.Range("B14").Value = "All"
$ws.Range("C14").Value = "This week"

$ws.Range("A15").Value = "Create focus group of family and friends that we can keep talking to as our project develops"
$ws.Range("B15").Value = "All"
$ws.Range("C15").Value = "Ongoing"

$ws.Range("A16").Value = "Become more familiar with the mean stack"
$ws.Range("B16").Value = "All"
$ws.Range("C16").Value = "This week"

# ---------------------------------------------------------------------
# 6. Decisions table (rows 19:21) - fix styles, then values.
# ---------------------------------------------------------------------
$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 30
$ws.Rows.Item(21).RowHeight = 45

$ws.Range("B19").Value = $null
$ws.Range("C19").Value = $null

$ws.Range("A19").Value = "Perhaps have the globe coming out of the book?"
$ws.Range("A20").Value = "Wenda to start working on understanding back end"
$ws.Range("A21").Value = "Hamza and Emily focussing on design, Stanni''s job to implement it"

# ---------------------------------------------------------------------
# 7. Sheet view: "1603" becomes the active / selected sheet.
# ---------------------------------------------------------------------
$ws.Select()
$ws.Range("A17:XFD17").Select()
$excel.ActiveWindow.ScrollRow = 6

$template.Select()
$template.Range("A1:C22").Select()
$excel.ActiveWindow.ScrollRow = 14

$ws.Select()
